# Adjust movetypes and movetype selection
$wb = $excel.ActiveWorkbook

# --- Update the "Maintenance" sheet's movetype values (column E / w3) ---
$ws = $wb.Worksheets.Item("Maintenance")

for ($row = 2; $row -le 79; $row++) {
    $d = $ws.Cells.Item($row, 4).Value()
    if ($d -eq 50) {
        $ws.Cells.Item($row, 5).Value = 200
    }
}

# --- Switch the active/selected sheet from Sheet2 to Maintenance ---
$ws.Activate()
$ws.Range("N6").Select()

$wb.Save()
